$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.392.70"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "1.566.63"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.17"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.500"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.96"
$ws.Range("E8").Value = "  -1.49%  "
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "1.789.60"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").Value = "1.578.99"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.83"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.516"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.44"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").Value = "27.390.07"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "212.93"
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.54"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.45"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.98"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.21"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").Value = "1.373.56"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("E35").Value = "  +1.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.967"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.532"
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.975"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.80"
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.07"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("E45").Value = "  +1.05%  "
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("D47").Value = "1.702.17"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.64"
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("D49").Value = "0.0₇0987"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0954"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("E51").Value = "  -0.79%  "
